$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 708.7
$ws.Range("I28").Value = 641.1429000000001
$ws.Range("J28").Value = 866.3333
$ws.Range("K28").Value = 641.1429000000001
$ws.Range("L28").Value = 866.3333
$ws.Range("M28").Value = -156.1429000000001
$ws.Range("N28").Value = -1836.3333
$ws.Range("H33").Value = 211.16667
$ws.Range("I33").Value = 173.4
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 173.4
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = 55.59999999999999
$ws.Range("N33").Value = -858
$ws.Range("H62").Value = 4474.5
$ws.Range("I62").Value = 802
$ws.Range("J62").Value = 4999.143
$ws.Range("K62").Value = 802
$ws.Range("L62").Value = 4999.143
$ws.Range("M62").Value = -178
$ws.Range("N62").Value = -6247.143
$ws.Range("H65").Value = 4474.5
$ws.Range("I65").Value = 802
$ws.Range("J65").Value = 4999.143
$ws.Range("K65").Value = 4010
$ws.Range("L65").Value = 24995.715
$ws.Range("M65").Value = -890
$ws.Range("N65").Value = -31235.715
$ws.Range("H106").Value = 27733.438
$ws.Range("I106").Value = 24582.334
$ws.Range("K106").Value = 24582.334
$ws.Range("M106").Value = -23951.334
$ws.Range("H137").Value = 2499.5
$ws.Range("I137").Value = 2499
$ws.Range("K137").Value = 7497
$ws.Range("M137").Value = -4947

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7752.778
$ws.Range("I61").Value = 8096.875
$ws.Range("K61").Value = 8096.875
$ws.Range("M61").Value = -7884.875
$ws.Range("H88").Value = 783.7857
$ws.Range("J88").Value = 1462.2
$ws.Range("L88").Value = 1462.2
$ws.Range("N88").Value = -2274.2
$ws.Range("H91").Value = 783.7857
$ws.Range("J91").Value = 1462.2
$ws.Range("L91").Value = 1462.2
$ws.Range("N91").Value = -4270.2
$ws.Range("H102").Value = 7411221.5
$ws.Range("I102").Value = 7940452
$ws.Range("K102").Value = 7940452
$ws.Range("M102").Value = -7938830
$ws.Range("H132").Value = 4724.8335
$ws.Range("J132").Value = 4945.5
$ws.Range("L132").Value = 14836.5
$ws.Range("N132").Value = -19896.5
$ws.Range("H135").Value = 149999.5
$ws.Range("J135").Value = 149999.5
$ws.Range("L135").Value = 149999.5
$ws.Range("N135").Value = -160139.5
$ws.Range("H136").Value = 7752.778
$ws.Range("I136").Value = 8096.875
$ws.Range("K136").Value = 24290.625
$ws.Range("M136").Value = -21740.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 703
$ws.Range("I99").Value = 585.0909
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 585.0909
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 912.9091
$ws.Range("N99").Value = -4996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 110079.164
$ws.Range("I22").Value = 168564.14
$ws.Range("J22").Value = 28200.2
$ws.Range("K22").Value = 168564.14
$ws.Range("L22").Value = 28200.2
$ws.Range("M22").Value = -168214.14
$ws.Range("N22").Value = -28900.2
$ws.Range("H29").Value = 28500
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 28500
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 28500
$ws.Range("M29").ClearContents() | Out-Null
$ws.Range("N29").Value = -29086
$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -27470
$ws.Range("N132").ClearContents() | Out-Null
$ws.Range("H134").Value = 2409.889
$ws.Range("I134").Value = 1801.2307
$ws.Range("K134").Value = 5403.6921
$ws.Range("M134").Value = -2868.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 250.24
$ws.Range("I5").Value = 250.24
$ws.Range("K5").Value = 750.72
$ws.Range("M5").Value = -638.72
$ws.Range("H10").Value = 52.857143
$ws.Range("I10").Value = 60
$ws.Range("J10").Value = 10
$ws.Range("K10").Value = 180
$ws.Range("L10").Value = 30
$ws.Range("M10").Value = -41
$ws.Range("N10").Value = -308
$ws.Range("H45").Value = 30
$ws.Range("I45").Value = 30
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 90
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 442
$ws.Range("N45").ClearContents() | Out-Null
$ws.Range("H69").Value = 3825
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622
$ws.Range("H72").Value = 3825
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents() | Out-Null
$ws.Range("H135").Value = 250.24
$ws.Range("I135").Value = 250.24
$ws.Range("K135").Value = 2252.16
$ws.Range("M135").Value = 282.8400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1090.5333
$ws.Range("I102").Value = 946.5833
$ws.Range("K102").Value = 946.5833
$ws.Range("M102").Value = 675.4167
$ws.Range("H132").Value = 2867.8147
$ws.Range("I132").Value = 3034.625
$ws.Range("J132").Value = 1533.3334
$ws.Range("K132").Value = 9103.875
$ws.Range("L132").Value = 4600.0002
$ws.Range("M132").Value = -6573.875
$ws.Range("N132").Value = -9660.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3839.3333
$ws.Range("I4").Value = 3839.3333
$ws.Range("K4").Value = 3839.3333
$ws.Range("M4").Value = -3726.3333
$ws.Range("H18").Value = 60000
$ws.Range("J18").Value = 60000
$ws.Range("L18").Value = 60000
$ws.Range("N18").Value = -60344
$ws.Range("H23").Value = 50019000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents() | Out-Null
$ws.Range("H28").Value = 3839.3333
$ws.Range("I28").Value = 3839.3333
$ws.Range("K28").Value = 3839.3333
$ws.Range("M28").Value = -3607.3333
$ws.Range("H33").Value = 15000000
$ws.Range("I33").Value = 15000000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 15000000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -14999710
$ws.Range("N33").ClearContents() | Out-Null
$ws.Range("H37").Value = 3839.3333
$ws.Range("I37").Value = 3839.3333
$ws.Range("K37").Value = 3839.3333
$ws.Range("M37").Value = -3732.3333
$ws.Range("H47").Value = 25200
$ws.Range("I47").Value = 21000
$ws.Range("K47").Value = 21000
$ws.Range("M47").Value = -20510
$ws.Range("H52").Value = 25200
$ws.Range("I52").Value = 21000
$ws.Range("K52").Value = 21000
$ws.Range("M52").Value = -20767

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents() | Out-Null
$ws.Range("N9").ClearContents() | Out-Null
$ws.Range("H11").Value = 17000
$ws.Range("J11").Value = 17000
$ws.Range("L11").Value = 17000
$ws.Range("N11").Value = -17284
$ws.Range("H12").Value = 14121
$ws.Range("I12").Value = 17995
$ws.Range("J12").Value = 10247
$ws.Range("K12").Value = 17995
$ws.Range("L12").Value = 10247
$ws.Range("M12").Value = -17853
$ws.Range("N12").Value = -10531
$ws.Range("H33").Value = 30607
$ws.Range("J33").Value = 31410.5
$ws.Range("L33").Value = 31410.5
$ws.Range("N33").Value = -31910.5
$ws.Range("H36").Value = 30607
$ws.Range("J36").Value = 31410.5
$ws.Range("L36").Value = 31410.5
$ws.Range("N36").Value = -31910.5
$ws.Range("H40").Value = 35000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 35000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 35000
$ws.Range("M40").ClearContents() | Out-Null
$ws.Range("N40").Value = -35298
$ws.Range("H47").Value = 18499
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 18499
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 18499
$ws.Range("M47").ClearContents() | Out-Null
$ws.Range("N47").Value = -19643
$ws.Range("H81").Value = 3475
$ws.Range("I81").Value = 3475
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6950
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5889
$ws.Range("N81").ClearContents() | Out-Null
$ws.Range("H84").Value = 3475
$ws.Range("I84").Value = 3475
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 34750
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -29446
$ws.Range("N84").ClearContents() | Out-Null
$ws.Range("H107").Value = 7969.4287
$ws.Range("I107").Value = 7969.4287
$ws.Range("K107").Value = 23908.2861
$ws.Range("M107").Value = -21988.2861
